$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19792.54129835656
$ws.Range("C2").Value = 15066.47555447738
$ws.Range("D2").Value = 621.6720529563064
$ws.Range("E2").Value = 30290.74036012588
$ws.Range("F2").Value = 5701.982570874461

$ws.Range("B3").Value = 11692.00206736152
$ws.Range("C3").Value = 7093.838029203962
$ws.Range("D3").Value = 302.4603624540867
$ws.Range("E3").Value = 16324.85947714854
$ws.Range("F3").Value = 3109.280991178463

$ws.Range("B4").Value = 6131.186826176758
$ws.Range("C4").Value = 7385.415563697593
$ws.Range("D4").Value = 309.6276001236665
$ws.Range("E4").Value = 13515.18985424745
$ws.Range("F4").Value = 2312.637594047707

$ws.Range("B5").Value = 946.1838048182954
$ws.Range("C5").Value = 545.0428615758075
$ws.Range("D5").Value = 9.291590378553938
$ws.Range("E5").Value = 447.9448287298518
$ws.Range("F5").Value = 280.063985648286
